$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 25 (old extra rows no longer needed)
$ws.Range("A7:A25").EntireRow.Delete()

# Update A2:A6 with the new combined tuple-style strings
$ws.Range("A2").Value = "('Demonic Tutor', ['{1}{B}', 'Sorcery', 'Search your library for a card, put that card into your hand, then shuffle your library.'])"
$ws.Range("A3").Value = "('Goblin Piledriver', ['{1}{R}', 'Creature " + [char]0x2014 + " Goblin Warrior', 'Protection from blue (This creature can" + [char]0x2019 + "t be blocked, targeted, dealt damage, or enchanted by anything blue.)', 'Whenever Goblin Piledriver attacks, it gets +2/+0 until end of turn for each other attacking Goblin.', '1/2'])"
$ws.Range("A4").Value = "('Living Wish', ['{1}{G}', 'Sorcery', 'You may reveal a creature or land card you own from outside the game and put it into your hand. Exile Living Wish.'])"
$ws.Range("A5").Value = "(`"Mind's Desire`", ['{4}{U}{U}', 'Sorcery', 'Shuffle your library. Then exile the top card of your library. Until end of turn, you may play that card without paying its mana cost. (If it has X in its mana cost, X is 0.)', 'Storm (When you cast this spell, copy it for each spell cast before it this turn.)'])"
$ws.Range("A6").Value = "(`"Orim's Chant`", ['{W}', 'Instant', 'Kicker {W} (You may pay an additional {W} as you cast this spell.)', 'Target player can" + [char]0x2019 + "t cast spells this turn. If this spell was kicked, creatures can" + [char]0x2019 + "t attack this turn.'])"
